$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 229, shifting existing rows 229:336 down to 230:337.
$ws.Rows.Item(229).Insert()

# Populate the newly inserted row 229 with the new record.
$ws.Cells.Item(229, 1).Value = 9
$ws.Cells.Item(229, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(229, 3).Value = "Metropolitana"
$ws.Cells.Item(229, 4).Value = 44636
$ws.Cells.Item(229, 5).Value = 13
$ws.Cells.Item(229, 6).Value = 100112044
$ws.Cells.Item(229, 7).Value = "Perejil"
$ws.Cells.Item(229, 8).Value = "Sin especificar"
$ws.Cells.Item(229, 9).Value = "Primera"
$ws.Cells.Item(229, 10).Value = 61
$ws.Cells.Item(229, 11).Value = 12000
$ws.Cells.Item(229, 12).Value = 14000
$ws.Cells.Item(229, 13).Value = 12984
$ws.Cells.Item(229, 14).Value = "$/docena de atados"
$ws.Cells.Item(229, 15).Value = "Región Metropolitana"
$ws.Cells.Item(229, 16).Value = 4328
$ws.Cells.Item(229, 17).Value = 3
$ws.Cells.Item(229, 18).Value = "Hortaliza"
